# feat: add 2022-Q4 data
#
# 1. Insert a brand-new worksheet named "2022-Q4" right after "总计" and
#    before "2022-Q3" (all the quarter sheets after it simply shift right
#    by one position, keeping their original name/content).
# 2. Populate "2022-Q4" with the new fund-holding table.
# 3. Update the "总计" (totals) sheet: insert a new top data row for
#    "2022-Q4" and shift the previous rows down by one (the row that falls
#    off the bottom, "2021-Q2", is re-appended as the new last row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the new sheet and move it into position 2 (right after
# "总计", i.e. right before the sheet that is currently "2022-Q3").
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q4"
$newSheet.Move($wb.Worksheets.Item(2))

# ---------------------------------------------------------------------
# Step 2: fill in the "2022-Q4" worksheet contents.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("2022-Q4")

# Header row (row 1) - same headers used by every other quarterly sheet.
$ws.Cells.Item(1, 2).Value = "基金代码"
$ws.Cells.Item(1, 3).Value = "基金名称"
$ws.Cells.Item(1, 4).Value = "基金规模"
$ws.Cells.Item(1, 5).Value = "股票总仓位"
$ws.Cells.Item(1, 6).Value = "仓位占比"
$ws.Cells.Item(1, 7).Value = "持有市值(亿元)"
$ws.Cells.Item(1, 8).Value = "仓位排名"
$ws.Range("B1:H1").Style = $wb.Worksheets.Item("2022-Q3").Range("B1:H1").Style

# Columns B, D, E, F, G hold numeric-looking text that must stay text
# (leading zeros / trailing zeros matter), so force a text format before
# assigning the values - otherwise Excel auto-coerces them to numbers.
$textCols = @(2, 4, 5, 6, 7)
foreach ($col in $textCols) {
    $ws.Range($ws.Cells.Item(2, $col), $ws.Cells.Item(5, $col)).NumberFormat = "@"
}

# Row 2
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "012868"
$ws.Cells.Item(2, 3).Value = "易方达标普信息科技指数（QDII-LOF）人民币 C"
$ws.Cells.Item(2, 4).Value = "5.09"
$ws.Cells.Item(2, 5).Value = "91.36"
$ws.Cells.Item(2, 6).Value = "2.16"
$ws.Cells.Item(2, 7).Value = "0.1099"
$ws.Cells.Item(2, 8).Value = 7

# Row 3
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "161128"
$ws.Cells.Item(3, 3).Value = "易方达标普信息科技指数（QDII-LOF）人民币"
$ws.Cells.Item(3, 4).Value = "5.09"
$ws.Cells.Item(3, 5).Value = "91.36"
$ws.Cells.Item(3, 6).Value = "2.16"
$ws.Cells.Item(3, 7).Value = "0.1099"
$ws.Cells.Item(3, 8).Value = 7

# Row 4
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = "003721"
$ws.Cells.Item(4, 3).Value = "易方达标普信息科技指数（QDII-LOF）美元A"
$ws.Cells.Item(4, 4).Value = "4.93"
$ws.Cells.Item(4, 5).Value = "91.36"
$ws.Cells.Item(4, 6).Value = "2.16"
$ws.Cells.Item(4, 7).Value = "0.1065"
$ws.Cells.Item(4, 8).Value = 7

# Row 5
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = "012869"
$ws.Cells.Item(5, 3).Value = "易方达标普信息科技指数（QDII-LOF）美元 C"
$ws.Cells.Item(5, 4).Value = "0.16"
$ws.Cells.Item(5, 5).Value = "91.36"
$ws.Cells.Item(5, 6).Value = "2.16"
$ws.Cells.Item(5, 7).Value = "0.0035"
$ws.Cells.Item(5, 8).Value = 7

$ws.Range("A2:A5").Style = $wb.Worksheets.Item("2022-Q3").Range("A2:A5").Style

# ---------------------------------------------------------------------
# Step 3: update the "总计" (totals) sheet - insert the new "2022-Q4" row
# at the top of the data and push every other row down by one. The row
# that used to be last ("2021-Q2") is re-written as the new row 8.
# ---------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")

# Existing data (top to bottom) before the edit:
#   row2: 2022-Q3 / 4 / 0.29
#   row3: 2022-Q2 / 4 / 0.2
#   row4: 2022-Q1 / 3 / 0.36
#   row5: 2021-Q4 / 7 / 4.73
#   row6: 2021-Q3 / 2 / 0.26
#   row7: 2021-Q2 / 2 / 0.32
# New data (top to bottom) after the edit:
#   row2: 2022-Q4 / 4 / 0.33   (new)
#   row3: 2022-Q3 / 4 / 0.29
#   row4: 2022-Q2 / 4 / 0.2
#   row5: 2022-Q1 / 3 / 0.36
#   row6: 2021-Q4 / 7 / 4.73
#   row7: 2021-Q3 / 2 / 0.26
#   row8: 2021-Q2 / 2 / 0.32   (new row, same values as old row7)

$totRows = @(
    @("2022-Q4", 4, 0.33),
    @("2022-Q3", 4, 0.29),
    @("2022-Q2", 4, 0.2),
    @("2022-Q1", 3, 0.36),
    @("2021-Q4", 7, 4.73),
    @("2021-Q3", 2, 0.26),
    @("2021-Q2", 2, 0.32)
)

for ($i = 0; $i -lt $totRows.Count; $i++) {
    $r = $i + 2
    $tot.Cells.Item($r, 1).Value = $i
    $tot.Cells.Item($r, 2).Value = $totRows[$i][0]
    $tot.Cells.Item($r, 3).Value = $totRows[$i][1]
    $tot.Cells.Item($r, 4).Value = $totRows[$i][2]
}

$tot.Range("A2:A8").Style = $tot.Range("A2").Style
